$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some Price-column values (e.g. "308.99") look like plain
# floating point numbers. Assigning them straight to .Value lets Excel
# auto-convert the cell to a Number (dropping significant trailing
# zeros / introducing binary rounding noise). Forcing the cell to Text
# format first makes Excel keep the literal string instead.
function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

$ws.Range('D2').Value = '45.870.19'
$ws.Range('E2').Value = '  -0.93%  '
$ws.Range('D3').Value = '2.604.99'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('E4').Value = '  -0.04%  '
Set-TextCell 'D5' '308.99'
$ws.Range('E5').Value = '  -1.36%  '
Set-TextCell 'D6' '98.79'
$ws.Range('E6').Value = '  -2.19%  '
$ws.Range('E7').Value = '  -1.01%  '
Set-TextCell 'D8' '0.999'
Set-TextCell 'D9' '0.577'
$ws.Range('E9').Value = '  -1.40%  '
Set-TextCell 'D10' '38.64'
$ws.Range('E10').Value = '  -0.65%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell 'D11' '0.0839'
$ws.Range('E11').Value = '  -0.49%  '
$ws.Range('B12').Value = 'OKB'
$ws.Range('C12').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell 'D12' '53.93'
$ws.Range('E12').Value = '  -0.50%  '
Set-TextCell 'D13' '8.05'
$ws.Range('E13').Value = '  -3.27%  '
$ws.Range('D14').Value = '2.997.53'
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('E15').Value = '  +0.81%  '
$ws.Range('D16').Value = '2.600.29'
$ws.Range('E16').Value = '  -0.46%  '
$ws.Range('E17').Value = '  -0.06%  '
Set-TextCell 'D18' '14.76'
$ws.Range('E18').Value = '  -2.40%  '
$ws.Range('D19').Value = '45.833.05'
$ws.Range('E19').Value = '  -1.45%  '
$ws.Range('E20').Value = '  -1.40%  '
Set-TextCell 'D21' '6.70'
$ws.Range('E21').Value = '  -1.09%  '
$ws.Range('E22').Value = '  -5.08%  '
Set-TextCell 'D23' '284.63'
$ws.Range('E23').Value = '  +11.58%  '
Set-TextCell 'D24' '73.84'
$ws.Range('E24').Value = '  +3.94%  '
Set-TextCell 'D25' '3.02'
$ws.Range('E25').Value = '  -1.96%  '
Set-TextCell 'D26' '2.25'
$ws.Range('E26').Value = '  +1.74%  '
Set-TextCell 'D27' '29.21'
$ws.Range('E27').Value = '  +3.41%  '
$ws.Range('E28').Value = '  +0.08%  '
Set-TextCell 'D29' '4.06'
$ws.Range('E29').Value = '  +0.42%  '
Set-TextCell 'D30' '10.54'
$ws.Range('E30').Value = '  -0.49%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell 'D31' '38.53'
$ws.Range('E31').Value = '  -5.30%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 'D32' '2.20'
$ws.Range('E32').Value = '  -3.19%  '
$ws.Range('E33').Value = '  +0.06%  '
Set-TextCell 'D34' '3.64'
$ws.Range('E34').Value = '  -2.16%  '
Set-TextCell 'D35' '157.27'
$ws.Range('E35').Value = '  +2.38%  '
Set-TextCell 'D36' '2.27'
$ws.Range('E36').Value = '  -0.77%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell 'D37' '2.80'
$ws.Range('E37').Value = '  -2.67%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 'D38' '0.0831'
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('E39').Value = '  +2.41%  '
$ws.Range('E40').Value = '  +0.11%  '
Set-TextCell 'D41' '15.96'
$ws.Range('E41').Value = '  -6.42%  '
$ws.Range('E42').Value = '  -0.76%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D43' '4.01'
$ws.Range('E43').Value = '  -5.32%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D44' '21.41'
$ws.Range('E44').Value = '  +0.38%  '
$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 'D45' '3.52'
$ws.Range('E45').Value = '  -3.04%  '
$ws.Range('D46').Value = '2.103.89'
$ws.Range('E46').Value = '  +3.22%  '
Set-TextCell 'D47' '0.998'
$ws.Range('E47').Value = '  -0.08%  '
Set-TextCell 'D48' '93.84'
$ws.Range('E48').Value = '  +2.54%  '
Set-TextCell 'D49' '9.18'
$ws.Range('E49').Value = '  -0.90%  '
Set-TextCell 'D50' '108.88'
$ws.Range('E50').Value = '  -2.83%  '
$ws.Range('D51').Value = '2.849.54'
$ws.Range('E51').Value = '  -0.56%  '
